$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 7885
$ws.Range("F5").Value = 7885
$ws.Range("F7").Value = 121
$ws.Range("F8").Value = 2132
$ws.Range("F9").Value = 8645
$ws.Range("F10").Value = 8645
$ws.Range("F13").Value = 91
$ws.Range("F14").Value = 5783
$ws.Range("F15").Value = 65
$ws.Range("F16").Value = 2776
$ws.Range("F17").Value = 1206
$ws.Range("F18").Value = 418
$ws.Range("F20").Value = 46
$ws.Range("F21").Value = 621
$ws.Range("F22").Value = 90
$ws.Range("F23").Value = 3945
$ws.Range("F26").Value = 63
$ws.Range("F28").Value = 183
$ws.Range("F30").Value = 5578
$ws.Range("F31").Value = 11
$ws.Range("F34").Value = 395
$ws.Range("F35").Value = 164
$ws.Range("F36").Value = 402
$ws.Range("F37").Value = 2713
$ws.Range("F40").Value = 1124
$ws.Range("F41").Value = 5153
$ws.Range("F42").Value = 79
$ws.Range("F44").Value = 47
$ws.Range("F45").Value = 3631
$ws.Range("F47").Value = 2348

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 17
$ws.Range("F10").Value = 132

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 271
$ws.Range("F3").Value = 1371

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 271
$ws.Range("F3").Value = 1371
$ws.Range("F5").Value = 7885
$ws.Range("F6").Value = 7885
$ws.Range("F7").Value = 121
$ws.Range("F8").Value = 2132
$ws.Range("F9").Value = 8645
$ws.Range("F10").Value = 8645
$ws.Range("F12").Value = 91
$ws.Range("F13").Value = 5783
$ws.Range("F14").Value = 65
$ws.Range("F15").Value = 2776
$ws.Range("F16").Value = 1206
$ws.Range("F17").Value = 418
$ws.Range("F19").Value = 46
$ws.Range("F21").Value = 621
$ws.Range("F22").Value = 90
$ws.Range("F23").Value = 3945
$ws.Range("F26").Value = 63
$ws.Range("F28").Value = 183
$ws.Range("F30").Value = 5578
$ws.Range("F31").Value = 11
$ws.Range("F33").Value = 395
$ws.Range("F34").Value = 164
$ws.Range("F35").Value = 402
$ws.Range("F37").Value = 2714
$ws.Range("F40").Value = 1124
$ws.Range("F42").Value = 5153
$ws.Range("F43").Value = 79
$ws.Range("F45").Value = 47
$ws.Range("F46").Value = 3631
$ws.Range("F47").Value = 2348
$ws.Range("F50").Value = 132
